$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 4941.6665
$ws.Range("I32").Value2 = 2975.5
$ws.Range("J32").Value2 = 5924.75
$ws.Range("K32").Value2 = 2975.5
$ws.Range("L32").Value2 = 5924.75
$ws.Range("M32").Value2 = -2649.5
$ws.Range("N32").Value2 = -6576.75
$ws.Range("H62").Value2 = 39499.4
$ws.Range("I62").Value2 = 20000
$ws.Range("J62").Value2 = 41666
$ws.Range("K62").Value2 = 20000
$ws.Range("L62").Value2 = 41666
$ws.Range("M62").Value2 = -19376
$ws.Range("N62").Value2 = -42914
$ws.Range("H65").Value2 = 39499.4
$ws.Range("I65").Value2 = 20000
$ws.Range("J65").Value2 = 41666
$ws.Range("K65").Value2 = 100000
$ws.Range("L65").Value2 = 208330
$ws.Range("M65").Value2 = -96880
$ws.Range("N65").Value2 = -214570
$ws.Range("H98").Value2 = 42361.5
$ws.Range("I98").Value2 = 46815.75
$ws.Range("K98").Value2 = 46815.75
$ws.Range("M98").Value2 = -45317.75
$ws.Range("H100").Value2 = 163949.83
$ws.Range("I100").Value2 = 100939.8
$ws.Range("J100").Value2 = 479000
$ws.Range("K100").Value2 = 100939.8
$ws.Range("L100").Value2 = 479000
$ws.Range("M100").Value2 = -100398.8
$ws.Range("N100").Value2 = -480082
$ws.Range("H112").Value2 = 68553.53
$ws.Range("J112").Value2 = 73353.78999999999
$ws.Range("L112").Value2 = 220061.37
$ws.Range("N112").Value2 = -222277.37
$ws.Range("H116").Value2 = 1015350.44
$ws.Range("I116").Value2 = 2781776.2
$ws.Range("K116").Value2 = 2781776.2
$ws.Range("M116").Value2 = -2778334.2
$ws.Range("H122").Value2 = 42361.5
$ws.Range("I122").Value2 = 46815.75
$ws.Range("K122").Value2 = 140447.25
$ws.Range("M122").Value2 = -137997.25
$ws.Range("H135").Value2 = 7474.3335
$ws.Range("I135").Value2 = 7474.3335
$ws.Range("K135").Value2 = 67269.0015
$ws.Range("M135").Value2 = -64734.0015
$ws.Range("H137").Value2 = 6387.171
$ws.Range("I137").Value2 = 8081.067
$ws.Range("J137").Value2 = 1767.4546
$ws.Range("K137").Value2 = 24243.201
$ws.Range("L137").Value2 = 5302.3638
$ws.Range("M137").Value2 = -21693.201
$ws.Range("N137").Value2 = -10402.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2056.0852
$ws.Range("I32").Value2 = 2089.913
$ws.Range("K32").Value2 = 2089.913
$ws.Range("M32").Value2 = -1802.913
$ws.Range("H45").Value2 = 6315
$ws.Range("I45").Value2 = 5614
$ws.Range("J45").Value2 = 7249.6665
$ws.Range("K45").Value2 = 5614
$ws.Range("L45").Value2 = 7249.6665
$ws.Range("M45").Value2 = -5237
$ws.Range("N45").Value2 = -8003.6665
$ws.Range("H74").Value2 = 3924.6099
$ws.Range("I74").Value2 = 2523.1936
$ws.Range("J74").Value2 = 8269
$ws.Range("K74").Value2 = 2523.1936
$ws.Range("L74").Value2 = 8269
$ws.Range("M74").Value2 = -1649.1936
$ws.Range("N74").Value2 = -10017
$ws.Range("H77").Value2 = 3924.6099
$ws.Range("I77").Value2 = 2523.1936
$ws.Range("J77").Value2 = 8269
$ws.Range("K77").Value2 = 12615.968
$ws.Range("L77").Value2 = 41345
$ws.Range("M77").Value2 = -8247.968000000001
$ws.Range("N77").Value2 = -50081
$ws.Range("H110").Value2 = 2781.7144
$ws.Range("I110").Value2 = 1904.9
$ws.Range("K110").Value2 = 1904.9
$ws.Range("M110").Value2 = 140.0999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value2 = 49362
$ws.Range("I26").Value2 = 49362
$ws.Range("K26").Value2 = 49362
$ws.Range("M26").Value2 = -49070
$ws.Range("H94").Value2 = 2789.1875
$ws.Range("I94").Value2 = 3241.2222
$ws.Range("K94").Value2 = 3241.2222
$ws.Range("M94").Value2 = -2790.2222
$ws.Range("H107").Value2 = 4280.8
$ws.Range("I107").Value2 = 4200.8887
$ws.Range("J107").Value2 = 5000
$ws.Range("K107").Value2 = 4200.8887
$ws.Range("L107").Value2 = 5000
$ws.Range("M107").Value2 = -2280.8887
$ws.Range("N107").Value2 = -8840
$ws.Range("H134").Value2 = 3758.5
$ws.Range("I134").Value2 = 2572.9
$ws.Range("K134").Value2 = 7718.700000000001
$ws.Range("M134").Value2 = -5183.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3208.5757
$ws.Range("I31").Value2 = 2522.611
$ws.Range("J31").Value2 = 4031.7334
$ws.Range("K31").Value2 = 2522.611
$ws.Range("L31").Value2 = 4031.7334
$ws.Range("M31").Value2 = -2227.611
$ws.Range("N31").Value2 = -4621.7334
$ws.Range("H34").Value2 = 3208.5757
$ws.Range("I34").Value2 = 2522.611
$ws.Range("J34").Value2 = 4031.7334
$ws.Range("K34").Value2 = 2522.611
$ws.Range("L34").Value2 = 4031.7334
$ws.Range("M34").Value2 = -2320.611
$ws.Range("N34").Value2 = -4435.7334
$ws.Range("H58").Value2 = 1560.3103
$ws.Range("I58").Value2 = 1332
$ws.Range("J58").Value2 = 1933.909
$ws.Range("K58").Value2 = 1332
$ws.Range("L58").Value2 = 1933.909
$ws.Range("M58").Value2 = -1129
$ws.Range("N58").Value2 = -2339.909
$ws.Range("H132").Value2 = 23772.158
$ws.Range("I132").Value2 = 1869.9166
$ws.Range("J132").Value2 = 61318.855
$ws.Range("K132").Value2 = 5609.7498
$ws.Range("L132").Value2 = 183956.565
$ws.Range("M132").Value2 = -3079.7498
$ws.Range("N132").Value2 = -189016.565
$ws.Range("H136").Value2 = 1560.3103
$ws.Range("I136").Value2 = 1332
$ws.Range("J136").Value2 = 1933.909
$ws.Range("K136").Value2 = 3996
$ws.Range("L136").Value2 = 5801.727000000001
$ws.Range("M136").Value2 = -1446
$ws.Range("N136").Value2 = -10901.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value2 = 68490
$ws.Range("J63").Value2 = 68490
$ws.Range("L63").Value2 = 68490
$ws.Range("N63").Value2 = -69862
$ws.Range("H66").Value2 = 68490
$ws.Range("J66").Value2 = 68490
$ws.Range("L66").Value2 = 205470
$ws.Range("N66").Value2 = -212334
$ws.Range("H132").Value2 = 3367.4517
$ws.Range("I132").Value2 = 3071.3809
$ws.Range("K132").Value2 = 9214.1427
$ws.Range("M132").Value2 = -6684.1427

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 123746
$ws.Range("I40").Value2 = 400000
$ws.Range("K40").Value2 = 400000
$ws.Range("M40").Value2 = -399864
$ws.Range("H61").Value2 = 27999.8
$ws.Range("I61").Value2 = 28333
$ws.Range("K61").Value2 = 28333
$ws.Range("M61").Value2 = -28131
$ws.Range("H82").Value2 = 2080.5264
$ws.Range("I82").Value2 = 1984.4
$ws.Range("J82").Value2 = 2441
$ws.Range("K82").Value2 = 1984.4
$ws.Range("L82").Value2 = 2441
$ws.Range("M82").Value2 = -1623.4
$ws.Range("N82").Value2 = -3163
$ws.Range("H85").Value2 = 2080.5264
$ws.Range("I85").Value2 = 1984.4
$ws.Range("J85").Value2 = 2441
$ws.Range("K85").Value2 = 1984.4
$ws.Range("L85").Value2 = 2441
$ws.Range("M85").Value2 = -736.4000000000001
$ws.Range("N85").Value2 = -4937
$ws.Range("H113").Value2 = 27999.8
$ws.Range("I113").Value2 = 28333
$ws.Range("K113").Value2 = 28333
$ws.Range("M113").Value2 = -26163
$ws.Range("H122").Value2 = 4895.8184
$ws.Range("I122").Value2 = 4837
$ws.Range("K122").Value2 = 14511
$ws.Range("M122").Value2 = -12061
$ws.Range("H132").Value2 = 289726.25
$ws.Range("I132").Value2 = 482725.16
$ws.Range("K132").Value2 = 1448175.48
$ws.Range("M132").Value2 = -1445645.48
$ws.Range("H136").Value2 = 10195.7
$ws.Range("I136").Value2 = 2988
$ws.Range("K136").Value2 = 8964
$ws.Range("M136").Value2 = -6414

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 48389.46
$ws.Range("I107").Value2 = 2642.182
$ws.Range("K107").Value2 = 7926.545999999999
$ws.Range("M107").Value2 = -6006.545999999999
$ws.Range("H113").Value2 = 2723.182
$ws.Range("I113").Value2 = 993.7143
$ws.Range("K113").Value2 = 2981.1429
$ws.Range("M113").Value2 = -811.1428999999998
$ws.Range("H122").Value2 = 29832.143
$ws.Range("I122").Value2 = 5123
$ws.Range("K122").Value2 = 15369
$ws.Range("M122").Value2 = -12919
$ws.Range("H126").Value2 = 29123.625
$ws.Range("I126").Value2 = 45665.89
$ws.Range("J126").Value2 = 7855
$ws.Range("K126").Value2 = 136997.67
$ws.Range("L126").Value2 = 23565
$ws.Range("M126").Value2 = -134527.67
$ws.Range("N126").Value2 = -28505
$ws.Range("H132").Value2 = 11689.464
$ws.Range("I132").Value2 = 13575.634
$ws.Range("K132").Value2 = 40726.902
$ws.Range("M132").Value2 = -38196.902
$ws.Range("H136").Value2 = 3050.1428
$ws.Range("I136").Value2 = 868.6667
$ws.Range("J136").Value2 = 4686.25
$ws.Range("K136").Value2 = 2606.0001
$ws.Range("L136").Value2 = 14058.75
$ws.Range("M136").Value2 = -56.0001000000002
$ws.Range("N136").Value2 = -19158.75
